$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 192 (weekly update: new week's
# cherry price data pushes the existing rows 192-213 down to 194-215).
$ws.Rows("192:193").Insert()

# --- New row 192: Lapins, Primera, Provincia de Curicó ---
$ws.Cells.Item(192, 1).Value = 11
$ws.Cells.Item(192, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(192, 3).Value = "Bíobío"
$ws.Cells.Item(192, 4).Value = 45267
$ws.Cells.Item(192, 5).Value = 8
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value = 100103
$ws.Cells.Item(192, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(192, 9).Value = 100103001
$ws.Cells.Item(192, 10).Value = "Cereza"
$ws.Cells.Item(192, 11).Value = "Lapins"
$ws.Cells.Item(192, 12).Value = "Primera"
$ws.Cells.Item(192, 13).Value = 150
$ws.Cells.Item(192, 14).Value = 9000
$ws.Cells.Item(192, 15).Value = 11000
$ws.Cells.Item(192, 16).Value = 9667
$ws.Cells.Item(192, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(192, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(192, 19).Value = 967
$ws.Cells.Item(192, 20).Value = 10

# --- New row 193: Royal Dawn, Primera, Provincia de Curicó ---
$ws.Cells.Item(193, 1).Value = 11
$ws.Cells.Item(193, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(193, 3).Value = "Bíobío"
$ws.Cells.Item(193, 4).Value = 45267
$ws.Cells.Item(193, 5).Value = 8
$ws.Cells.Item(193, 6).Value = "Fruta"
$ws.Cells.Item(193, 7).Value = 100103
$ws.Cells.Item(193, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(193, 9).Value = 100103001
$ws.Cells.Item(193, 10).Value = "Cereza"
$ws.Cells.Item(193, 11).Value = "Royal Dawn"
$ws.Cells.Item(193, 12).Value = "Primera"
$ws.Cells.Item(193, 13).Value = 120
$ws.Cells.Item(193, 14).Value = 1000
$ws.Cells.Item(193, 15).Value = 1000
$ws.Cells.Item(193, 16).Value = 1000
$ws.Cells.Item(193, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(193, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(193, 19).Value = 100
$ws.Cells.Item(193, 20).Value = 10
